# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" before the current "2021-Q4" sheet
#    (position 2, right after "总计") and populate it with the fund
#    holding data for that quarter.
# 2. Prepend a corresponding summary row to the "总计" (totals) sheet,
#    shifting the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q3" worksheet right before "2021-Q4"
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q4Sheet)
$newSheet.Name = "2022-Q3"

# header row text
$newSheet.Range("B1:H1").NumberFormat = "@"
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# data rows text
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "160921"
$newSheet.Range("C2").Value = "大成多策略混合（LOF）A"
$newSheet.Range("D2").Value = "2.90"
$newSheet.Range("E2").Value = "87.66"
$newSheet.Range("F2").Value = "3.53"
$newSheet.Range("G2").Value = "0.1024"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "016062"
$newSheet.Range("C3").Value = "大成多策略混合（LOF）C"
$newSheet.Range("D3").Value = "1.37"
$newSheet.Range("E3").Value = "87.66"
$newSheet.Range("F3").Value = "3.53"
$newSheet.Range("G3").Value = "0.0484"
$newSheet.Range("H3").Value = 8

# Re-apply the same cell formatting the other quarter sheets use (header
# row bold+bordered style, "A" index column style) by copying it across,
# and clear the transient quote-prefix/text-numfmt formatting that typing
# the values above picked up on the data cells. Re-resolve the "2021-Q4"
# sheet by name (its position shifted when the new sheet was inserted, so
# the earlier $q4Sheet handle now points at "2022-Q3" itself).
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q4Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

$q4Sheet.Range("C2").Copy()
$newSheet.Range("B2:G3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: shift the "总计" rows down and insert the 2022-Q3 summary row
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q1"
$totalSheet.Range("C5").Value = 5
$totalSheet.Range("D5").Value = 1.62

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.02

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.03

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.15

# carry the "A" column style down onto the newly written cells
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# restore the original active sheet ("2021-Q1", the last tab) — adding a
# new sheet shifts Excel's "active sheet" focus onto it by default.
$wb.Worksheets.Item("2021-Q1").Activate()
